# This edit inserts a new data row (row 22) into the "Berenjena" sheet,
# pushing all existing rows 22:91 down to 23:92, and populates the new
# row 22 with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22 (shifts rows 22:91 -> 23:92)
$ws.Rows("22:22").Insert()

# Populate the new row 22 with the new observation
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44715
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112001
$ws.Range("G22").Value = "Berenjena"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 6500
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = 6700
$ws.Range("N22").Value = "`$/caja 60 unidades"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 112
$ws.Range("Q22").Value = 60
$ws.Range("R22").Value = "Hortaliza"
